$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / numeric values that Excel's type-sniffer leaves alone.
$ws.Range("A2").Value = "W0001"
$ws.Range("B2").Value = "VP "

# C2/D2 need literal text "2021-11-12" / "2021-11-26" (not Excel date
# serials). Writing the text straight into .Value lets Excel "helpfully"
# reinterpret it as a date and stamp a new number-format style onto the
# cell. Route it through a scratch text-formula cell instead: the
# formula result is a string, copying + pasting *values only* keeps it
# a plain string with the workbook's default style, and clearing the
# scratch cell afterwards leaves no residue.
$scratch = $ws.Range("Z1")

$scratch.Formula = "=""2021-11-12"""
$scratch.Copy()
$ws.Range("C2").PasteSpecial(-4163)

$scratch.Formula = "=""2021-11-26"""
$scratch.Copy()
$ws.Range("D2").PasteSpecial(-4163)

$scratch.ClearContents()

$ws.Range("E2").Value = 5
$ws.Range("F2").Value = "Instalasi jaringan"
$ws.Range("G2").Value = "Membuat Instalasi jaringan"

# Rows 3 and 4 (the other two WBS entries) are dropped entirely.
$ws.Rows("3:4").Delete()
